$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing row 1 (the lone "240404930" cell) down to row 2 by
# inserting a fresh row above it - this preserves the original cell's
# text/value/type exactly (no re-entry / auto-conversion side effects).
$ws.Rows(1).Insert()

$desc8200 = "- 96cpm (B&W) Includes: 200-sheet ADF, duplex unit, 3,000 sheets of paper capacity, and 2 x 250GB HDD. Toner & Developer not included with mainframe and must be added to initial order. Max Monthly Vol 1M. Requires Surge Protector (241007138MIU) & Pro Toner 8200S(244828484), as well as one of the following Finoshers SR5050 (241404548) or SR5060 (241404550) or SK5030 Stacker (241404656)."
$desc8210 = "- 111cpm (B&W) Includes: 200-sheet ADF, duplex unit, 3,000 sheets of paper capacity, and 2 x 250GB HDD. Toner & Developer not included with mainframe and must be added to initial order. Max Monthly Vol 1M. Requires Surge Protector (241007138MIU) & Pro Toner 8200S(244828484), as well as one of the following Finoshers SR5050 (241404548) or SR5060 (241404550) or SK5030 Stacker (241404656)."

# Row 1 - Pro 8200S.
# A1 carries a trailing newline and looks numeric, so force text entry via a
# Text number format (otherwise it is auto-parsed as a number), then reset
# the style back to the default "Normal" cell style so no stray formatting
# is left behind.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "240404929`n"
$ws.Range("A1").Style = "Normal"

$ws.Range("B1").Value = "Pro 8200S"
$ws.Range("C1").Value = $desc8200

$ws.Range("D1:G1").NumberFormat = "@"
$ws.Range("D1").Value = "26,488"
$ws.Range("E1").Value = "23,839"
$ws.Range("F1").Value = "22,514"
$ws.Range("G1").Value = "48,011"
$ws.Range("D1:G1").Style = "Normal"

# Row 2 - Pro 8210S (A2 already holds the original value after the insert).
$ws.Range("B2").Value = "Pro 8210S"
$ws.Range("C2").Value = $desc8210

$ws.Range("D2:G2").NumberFormat = "@"
$ws.Range("D2").Value = "33,285"
$ws.Range("E2").Value = "29,957"
$ws.Range("F2").Value = "28,293"
$ws.Range("G2").Value = "s59,880 "
$ws.Range("D2:G2").Style = "Normal"
